# eventbuttons.xlsx — add a new "alarm(n,<bool>)" command row to the
# Commands sheet, right after the "Artisan Command" / "$" row.
#
# Effect:
#   - a new row is inserted after row 95 (becomes row 96) with:
#       B96 = "alarm(n,<bool>)"
#       C96 = "enables/disables alarm number " & italic "n"
#   - every row from the old row 96 onward shifts down by one
#   - the sheet view scroll position / selection is updated

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row right before the current row 96 ("alarms(<bool>)"),
# i.e. directly after the "Artisan Command" block (row 95).
$ws.Rows.Item(96).Insert()

# Column B: function signature, italic style (matches the rest of column B)
$bCell = $ws.Cells.Item(96, 2)
$bCell.Value = "alarm(n,<bool>)"
$bCell.Font.Name = "Calibri"
$bCell.Font.Size = 11
$bCell.Font.Color = 0
$bCell.Font.Italic = $true

# Column C: description, regular style, with the trailing "n" placeholder
# in italics (matches similar two-run descriptions elsewhere, e.g. the
# "ramp(n,<bool>)" row).
$cCell = $ws.Cells.Item(96, 3)
$text = "enables/disables alarm number n"
$cCell.Value = $text
$cCell.Font.Name = "Calibri"
$cCell.Font.Size = 11
$cCell.Font.Color = 0
$cCell.Font.Italic = $false

$nPos = $text.Length  # 1-based index of the final "n" character
$cCell.Characters($nPos, 1).Font.Name = "Calibri"
$cCell.Characters($nPos, 1).Font.Size = 11
$cCell.Characters($nPos, 1).Font.Color = 0
$cCell.Characters($nPos, 1).Font.Italic = $true

# Update the visible scroll/selection state to match the edited document.
$ws.Application.GoTo($ws.Range("C98"), $true)
$ws.Range("C98").Select()
